# Update settings and data on the "gen" sheet.
#
# The three existing generator rows (4: nuclear, 5: NGCC, 6: NGCT) are
# re-ordered (NGCC -> row4, NGCT -> row5, nuclear -> row6) and the nuclear
# unit's EXUNITS/CAP values change from 5/1000 to 17/300. A new, empty row 7
# is also added below the data, and the sheet view's scroll/selection state
# is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gen")
$ws.Activate()

# --- Row 4: 1_ngcc_1 / Tech6 / NGCC / Gas CC / CC / Gas / 280 / 200 ---
$ws.Range("A4").Value = "1_ngcc_1"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Tech6"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "NGCC"
$ws.Range("F4").Value = "Gas CC"
$ws.Range("G4").Value = "CC"
$ws.Range("H4").Value = "Gas"
$ws.Range("I4").Value = 280
$ws.Range("J4").Value = 200

# --- Row 5: 1_ngct_1 / Tech7 / NGCT / Gas CT / CT / Gas / 508 / 50 ---
$ws.Range("A5").Value = "1_ngct_1"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "Tech7"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "NGCT"
$ws.Range("F5").Value = "Gas CT"
$ws.Range("G5").Value = "CT"
$ws.Range("H5").Value = "Gas"
$ws.Range("I5").Value = 508
$ws.Range("J5").Value = 50

# --- Row 6: 1_nuclear_1 / Tech10 / Nuclear_adv / Nuclear / Nuclear / Nuclear / 17 / 300 ---
$ws.Range("A6").Value = "1_nuclear_1"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Tech10"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "Nuclear_adv"
$ws.Range("F6").Value = "Nuclear"
$ws.Range("G6").Value = "Nuclear"
$ws.Range("H6").Value = "Nuclear"
$ws.Range("I6").Value = 17
$ws.Range("J6").Value = 300

# --- New trailing blank row so the used range grows to A1:BD7 ---
$ws.Range("A7").RowHeight = $ws.Range("A7").RowHeight

# --- Sheet view: scroll so column B is leftmost, select B6 ---
$ws.Range("B6").Select()
$excel.ActiveWindow.ScrollColumn = 2
